$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Means"
# Add two new columns (F = "Within 5 miles...", G = "Within 10 miles...")
# and update several existing values that changed alongside the new columns.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Means")

# New header cells
$ws1.Range("F1").Value = "Within 5 miles of HFC production facility"
$ws1.Range("G1").Value = "Within 10 miles of HFC production facility"

# Row 2 - % White
$ws1.Range("F2").Value = 92
$ws1.Range("G2").Value = 91

# Row 3 - % Black or African American
$ws1.Range("F3").Value = 2.2
$ws1.Range("G3").Value = 2.1

# Row 4 - % Other
$ws1.Range("F4").Value = 6.2
$ws1.Range("G4").Value = 7.1

# Row 5 - % Hispanic
$ws1.Range("F5").Value = 44
$ws1.Range("G5").Value = 40

# Row 6 - Median Income [1,000 2019$]
$ws1.Range("F6").Value = 69
$ws1.Range("G6").Value = 61

# Row 7 - % Below Poverty Line
$ws1.Range("F7").Value = 3.4
$ws1.Range("G7").Value = 6

# Row 8 - % Below Half the Poverty Line
$ws1.Range("F8").Value = 3.7
$ws1.Range("G8").Value = 4.9

# Row 9 - Total Cancer Risk (per million) -- values for B:E also changed
$ws1.Range("B9").Value = 29
$ws1.Range("C9").Value = 31
$ws1.Range("D9").Value = 20
$ws1.Range("E9").Value = 20
$ws1.Range("F9").Value = 20
$ws1.Range("G9").Value = 20

# Row 10 - Total Respiratory (hazard quotient) -- values for B:E also changed
$ws1.Range("B10").Value = 0.37
$ws1.Range("C10").Value = 0.36
$ws1.Range("D10").Value = 0.2
$ws1.Range("E10").Value = 0.2
$ws1.Range("F10").Value = 0.21
$ws1.Range("G10").Value = 0.21

# ---------------------------------------------------------------------------
# Sheet 2: "Standard Deviations"
# Same shape of change: two new columns F, G plus a handful of value updates.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Standard Deviations")

# New header cells
$ws2.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$ws2.Range("G1").Value = "Within 10 mile of HFC production facility SD"

# Row 2 - % White
$ws2.Range("F2").Value = 5.1
$ws2.Range("G2").Value = 8.3

# Row 3 - % Black or African American
$ws2.Range("F3").Value = 2.2
$ws2.Range("G3").Value = 2.9

# Row 4 - % Other
$ws2.Range("F4").Value = 3.9
$ws2.Range("G4").Value = 8.3

# Row 5 - % Hispanic
$ws2.Range("F5").Value = 20
$ws2.Range("G5").Value = 22

# Row 6 - Median Income [1,000 2019$]
$ws2.Range("F6").Value = 17
$ws2.Range("G6").Value = 16

# Row 7 - % Below Poverty Line
$ws2.Range("F7").Value = 3.8
$ws2.Range("G7").Value = 8.5

# Row 8 - % Below Half the Poverty Line
$ws2.Range("F8").Value = 3.7
$ws2.Range("G8").Value = 8.9

# Row 9 - Total Cancer Risk (per million) -- values for B:E also changed
$ws2.Range("B9").Value = 10
$ws2.Range("C9").Value = 14
$ws2.Range("D9").Value = 0
$ws2.Range("E9").Value = 0
$ws2.Range("F9").Value = 0
$ws2.Range("G9").Value = 0

# Row 10 - Total Respiratory (hazard quotient) -- values for B:E also changed
$ws2.Range("B10").Value = 0.14
$ws2.Range("C10").Value = 0.079
$ws2.Range("D10").Value = 0.00000000000000002
$ws2.Range("E10").Value = 0.000000000000000013
$ws2.Range("F10").Value = 0.03
$ws2.Range("G10").Value = 0.033
